$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the serial numbers that currently live in D2 before we restructure columns
$serials = $ws.Range("D2").Value2

# Drop the old per-card-type columns (D:I) - column C becomes the single summary column
$ws.Range("D1:I1").EntireColumn.Delete()

# Rename the remaining header and write the combined, labeled serial-number text
$ws.Range("C1").Value = "Assigned Serial Numbers"
$ws.Range("C2").Value = "Choices `$200: " + ($serials -replace " ", ", ")
$ws.Range("C2").NumberFormat = "@"

# Widen column C so the long serial-number list is easy to read
# (COM ColumnWidth is in "characters"; Excel stores it back out with a
#  ~5/6 character padding offset, so compensate to land on a stored width of 80)
$ws.Range("C:C").ColumnWidth = 79.16666666666667

"done"
